$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.618.90"
$ws.Range("E2").Value = "  +2.70%  "

$ws.Range("D3").Value = "1.855.28"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.75%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5239"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3296"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06752"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7747"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07702"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "

$ws.Range("D13").Value = "1.825.31"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.062"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("E16").Value = "  +0.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007917"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").Value = "26.702.25"
$ws.Range("E20").Value = "  +2.87%  "

$ws.Range("D21").Value = "2.073.16"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.618"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.754"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.365"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.663"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.226"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.215"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08779"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04881"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.144"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.890"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.157"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01824"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.252"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4977"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "115.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9115"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.077"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.35%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.784"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4308"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1298"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.198"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05927"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
